# Auto-generated Excel COM-interop script
# Applies updated market-price derived values (H..N columns) per row,
# as produced by the scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 21635
$ws.Range("J3").Value = 21635
$ws.Range("L3").Value = 21635
$ws.Range("N3").Value = -21863

$ws.Range("H39").Value = 356.9091
$ws.Range("I39").Value = 158.44444
$ws.Range("J39").Value = 1250
$ws.Range("K39").Value = 475.33332
$ws.Range("L39").Value = 3750
$ws.Range("M39").Value = -179.33332
$ws.Range("N39").Value = -4342

$ws.Range("H42").Value = 211.33333
$ws.Range("I42").Value = 257.5
$ws.Range("K42").Value = 772.5
$ws.Range("M42").Value = -542.5

$ws.Range("H100").Value = 3977.7778
$ws.Range("I100").Value = 3216.6667
$ws.Range("K100").Value = 3216.6667
$ws.Range("M100").Value = -2675.6667

$ws.Range("H102").Value = 21635
$ws.Range("J102").Value = 21635
$ws.Range("L102").Value = 21635
$ws.Range("N102").Value = -28125

$ws.Range("H132").Value = 1021.0417
$ws.Range("I132").Value = 1021.0417
$ws.Range("K132").Value = 3063.1251
$ws.Range("M132").Value = -533.1251000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4762455.5
$ws.Range("I32").Value = 578.4
$ws.Range("K32").Value = 578.4
$ws.Range("M32").Value = -291.4

$ws.Range("H61").Value = 4999.5
$ws.Range("I61").Value = 4999
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 4999
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -4787
$ws.Range("N61").Value = -5424

$ws.Range("H74").Value = 2902.8
$ws.Range("I74").Value = 2325.889
$ws.Range("K74").Value = 2325.889
$ws.Range("M74").Value = -1451.889

$ws.Range("H77").Value = 2902.8
$ws.Range("I77").Value = 2325.889
$ws.Range("K77").Value = 11629.445
$ws.Range("M77").Value = -7261.445

$ws.Range("H102").Value = 3843.0435
$ws.Range("I102").Value = 2311.875
$ws.Range("K102").Value = 2311.875
$ws.Range("M102").Value = -689.875

$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws.Range("H136").Value = 4999.5
$ws.Range("I136").Value = 4999
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 14997
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -12447
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 925.2
$ws.Range("I29").Value = 925.2
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 925.2
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -636.2
$ws.Range("N29").ClearContents()

$ws.Range("H36").Value = 933.125
$ws.Range("I36").Value = 933.125
$ws.Range("K36").Value = 933.125
$ws.Range("M36").Value = -399.125

$ws.Range("H96").Value = 19047
$ws.Range("I96").Value = 19047
$ws.Range("K96").Value = 19047
$ws.Range("M96").Value = -16301

$ws.Range("H103").Value = 15187.5
$ws.Range("J103").Value = 15187.5
$ws.Range("L103").Value = 15187.5
$ws.Range("N103").Value = -17531.5

$ws.Range("H105").Value = 2860716.5
$ws.Range("I105").Value = 2860716.5
$ws.Range("K105").Value = 2860716.5
$ws.Range("M105").Value = -2858969.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 3948.111
$ws.Range("J94").Value = 4976.2856
$ws.Range("L94").Value = 4976.2856
$ws.Range("N94").Value = -5878.2856

$ws.Range("H134").Value = 5696.8
$ws.Range("I134").Value = 4622.25
$ws.Range("K134").Value = 13866.75
$ws.Range("M134").Value = -11331.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 93.375
$ws.Range("J40").Value = 299.5
$ws.Range("L40").Value = 1198
$ws.Range("N40").Value = -1336

$ws.Range("H46").Value = 4000
$ws.Range("J46").Value = 4000
$ws.Range("L46").Value = 12000
$ws.Range("N46").Value = -12182

$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws.Range("H55").Value = 2214
$ws.Range("I55").Value = 397.5
$ws.Range("J55").Value = 3021.3333
$ws.Range("K55").Value = 1192.5
$ws.Range("L55").Value = 9063.999899999999
$ws.Range("M55").Value = -1015.5
$ws.Range("N55").Value = -9417.999899999999

$ws.Range("H92").Value = 2469
$ws.Range("I92").Value = 942
$ws.Range("J92").Value = 3487
$ws.Range("K92").Value = 2826
$ws.Range("L92").Value = 10461
$ws.Range("M92").Value = -1578
$ws.Range("N92").Value = -12957

$ws.Range("H114").Value = 1511.5454
$ws.Range("I114").Value = 1529
$ws.Range("J114").Value = 1497
$ws.Range("K114").Value = 4587
$ws.Range("L114").Value = 4491
$ws.Range("M114").Value = -1333
$ws.Range("N114").Value = -10999

$ws.Range("H131").Value = 1899.4286
$ws.Range("J131").Value = 2458
$ws.Range("L131").Value = 7374
$ws.Range("N131").Value = -17454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 135.8108
$ws.Range("J2").Value = 502.75
$ws.Range("L2").Value = 502.75
$ws.Range("N2").Value = -728.75

$ws.Range("H39").Value = 30000
$ws.Range("J39").Value = 30000
$ws.Range("L39").Value = 30000
$ws.Range("N39").Value = -31064

$ws.Range("H122").Value = 3409
$ws.Range("I122").Value = 3246
$ws.Range("J122").Value = 3653.5
$ws.Range("K122").Value = 9738
$ws.Range("L122").Value = 10960.5
$ws.Range("M122").Value = -7288
$ws.Range("N122").Value = -15860.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2868.2856
$ws.Range("I16").Value = 2879.8333
$ws.Range("J16").Value = 2799
$ws.Range("K16").Value = 2879.8333
$ws.Range("L16").Value = 2799
$ws.Range("M16").Value = -2709.8333
$ws.Range("N16").Value = -3139

$ws.Range("H55").Value = 1077.4375
$ws.Range("I55").Value = 1097.1111
$ws.Range("K55").Value = 1097.1111
$ws.Range("M55").Value = -924.1111000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 775.6667
$ws.Range("I113").Value = 543.7
$ws.Range("K113").Value = 1631.1
$ws.Range("M113").Value = 538.8999999999999

$ws.Range("H132").Value = 1633.3334
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470
